# Insert a new data row before row 620 (shifts existing rows 620-705 down
# to 621-706) and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(620).Insert()

$ws.Range("A620").Value = 5
$ws.Range("B620").Value = "Macroferia Regional de Talca"
$ws.Range("C620").Value = "Maule"
$ws.Range("D620").Value = 45131
$ws.Range("E620").Value = 7
$ws.Range("F620").Value = 100112043
$ws.Range("G620").Value = "Pepino ensalada"
$ws.Range("H620").Value = "Sin especificar"
$ws.Range("I620").Value = "Primera"
$ws.Range("J620").Value = 500
$ws.Range("K620").Value = 8000
$ws.Range("L620").Value = 8000
$ws.Range("M620").Value = 8000
$ws.Range("N620").Value = "`$/caja 60 unidades"
$ws.Range("O620").Value = "Región de Arica y Parinacota"
$ws.Range("P620").Value = 133
$ws.Range("Q620").Value = 60
$ws.Range("R620").Value = "Hortaliza"
